$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Clear the "Approved/Rejected" (I) and "ReasonToReject" (J) data columns
# for rows 2 through 30, leaving the header row (row 1) intact.
$ws.Range("I2:J30").ClearContents()

# Update the active selection to match the saved view state.
$ws.Range("M11").Select()
